$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds "Price" values that look numeric (e.g. "568.11") or
# pseudo-numeric with thousands separators (e.g. "63.137.56"). Excel's
# Range.Value setter auto-coerces plain decimal-looking strings into
# real numbers (losing the original text formatting / introducing float
# rounding noise), so force the cell to Text format first, assign the
# literal string, then drop back to the Normal style so no stray
# number-format override is left behind on the cell.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '63.137.56'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.548.11'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +3.09%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '568.11'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '146.91'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.25%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -0.36%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.545.44'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +3.04%  '
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('E11').Value = '  -2.01%  '
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('E13').Value = '  +0.37%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '27.52'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +3.39%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.001.54'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +3.00%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '63.032.09'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('E17').Value = '  +1.59%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.553.30'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +3.25%  '
$ws.Range('E19').Value = '  +1.51%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '335.37'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.62%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.33'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.72%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.78'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.63%  '
$ws.Range('E23').Value = '  -0.19%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '65.24'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.62%  '
$ws.Range('E25').Value = '  +8.97%  '
$ws.Range('E26').Value = '  -1.97%  '
$ws.Range('E27').Value = '  +7.74%  '
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('E29').Value = '  +4.02%  '
$ws.Range('E30').Value = '  +7.69%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0₃0820'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.73%  '
$ws.Range('E32').Value = '  +0.80%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '176.06'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('E34').Value = '  +3.49%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '413.55'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +11.06%  '
$ws.Range('E36').Value = '  +0.32%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '18.94'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.66%  '
$ws.Range('E38').Value = '  +0.16%  '
$ws.Range('E40').Value = '  +3.41%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '39.30'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.87%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '152.67'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.29%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '21.05'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +2.37%  '
$ws.Range('E46').Value = '  +0.56%  '
$ws.Range('E48').Value = '  +1.74%  '
$ws.Range('E49').Value = '  +5.06%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '18.35'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.13%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.78'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.12%  '
